# Applies the "cryptos list" price/volume refresh described in the commit
# message. Only the D (Price) and E (Volume(1h)) text columns change, for the
# 50 data rows (2-51). All values in this sheet are stored as text, including
# numeric-looking ones (e.g. "1.010"), so a plain Range.Value assignment is not
# safe for D-column values that Excel would auto-parse as a number (it would
# coerce type to Number and could drop significant trailing zeros, e.g.
# "1.010" -> 1.01). To keep these as literal text we:
#   1. assign with a leading apostrophe to force text entry, then
#   2. paste-format (xlPasteFormats = -4122) from a neighbouring untouched
#      cell on the same row (style index 0) back onto the cell, which clears
#      the "number stored as text" formatting Excel applies in step 1 while
#      leaving the literal text value untouched.
# The E-column values (e.g. "  +0.48%  ") already contain spaces/percent signs
# so Excel never tries to parse them as numbers; a direct Value assignment is
# sufficient there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.606.87"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Value = "  +0.48%  "

$ws.Range("D3").Value = "'2.116.79"
$ws.Range("B3").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Value = "  +1.18%  "

$ws.Range("D4").Value = "'1.010"
$ws.Range("B4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$ws.Range("E4").Value = "  +0.79%  "

$ws.Range("D5").Value = "'337.25"
$ws.Range("B5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").Value = "  +1.98%  "

$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("D7").Value = "'0.5247"
$ws.Range("B7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$ws.Range("E7").Value = "  +0.61%  "

$ws.Range("D8").Value = "'0.4566"
$ws.Range("B8").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").Value = "  +4.04%  "

$ws.Range("D9").Value = "'54.58"
$ws.Range("B9").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null
$ws.Range("E9").Value = "  +1.14%  "

$ws.Range("D10").Value = "'0.09147"
$ws.Range("B10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").Value = "  +2.43%  "

$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("D12").Value = "'24.57"
$ws.Range("B12").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").Value = "  +1.21%  "

$ws.Range("D13").Value = "'2.119.65"
$ws.Range("B13").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("E13").Value = "  +1.28%  "

$ws.Range("D14").Value = "'6.870"
$ws.Range("B14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null

$ws.Range("D15").Value = "'8.130"
$ws.Range("B15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = "  +5.89%  "

$ws.Range("D16").Value = "'0.00001178"
$ws.Range("B16").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Value = "  +4.89%  "

$ws.Range("D17").Value = "'97.23"
$ws.Range("B17").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").Value = "  +1.34%  "

$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("D19").Value = "'0.06692"
$ws.Range("B19").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null
$ws.Range("E19").Value = "  +1.42%  "

$ws.Range("D20").Value = "'19.43"
$ws.Range("B20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = "  +1.28%  "

$ws.Range("D21").Value = "'1.009"
$ws.Range("B21").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4122) | Out-Null
$ws.Range("E21").Value = "  +0.79%  "

$ws.Range("D22").Value = "'6.316"
$ws.Range("B22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = "  +0.83%  "

$ws.Range("D23").Value = "'30.657.22"
$ws.Range("B23").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = "  +0.52%  "

$ws.Range("D24").Value = "'12.87"
$ws.Range("B24").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4122) | Out-Null
$ws.Range("E24").Value = "  +4.70%  "

$ws.Range("D25").Value = "'2.363"
$ws.Range("B25").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null
$ws.Range("E25").Value = "  +1.90%  "

$ws.Range("D26").Value = "'2.365.75"
$ws.Range("B26").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = "  +1.23%  "

$ws.Range("D27").Value = "'22.43"
$ws.Range("B27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = "  +0.80%  "

$ws.Range("D28").Value = "'164.06"
$ws.Range("B28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("D29").Value = "'2.551"
$ws.Range("B29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = "  -0.29%  "

$ws.Range("D30").Value = "'134.24"
$ws.Range("B30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").Value = "  +2.04%  "

$ws.Range("D31").Value = "'1.218"
$ws.Range("B31").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4122) | Out-Null
$ws.Range("E31").Value = "  +2.70%  "

$ws.Range("D32").Value = "'0.1073"
$ws.Range("B32").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4122) | Out-Null

$ws.Range("D33").Value = "'1.654"
$ws.Range("B33").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4122) | Out-Null
$ws.Range("E33").Value = "  -0.32%  "

$ws.Range("D34").Value = "'6.380"
$ws.Range("B34").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4122) | Out-Null
$ws.Range("E34").Value = "  +3.47%  "

$ws.Range("D35").Value = "'3.956"
$ws.Range("B35").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4122) | Out-Null
$ws.Range("E35").Value = "  +1.47%  "

$ws.Range("E36").Value = "  +5.74%  "

$ws.Range("D37").Value = "'5.882"
$ws.Range("B37").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4122) | Out-Null
$ws.Range("E37").Value = "  +7.29%  "

$ws.Range("D38").Value = "'0.02630"
$ws.Range("B38").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4122) | Out-Null
$ws.Range("E38").Value = "  +2.31%  "

$ws.Range("D39").Value = "'0.06855"
$ws.Range("B39").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4122) | Out-Null
$ws.Range("E39").Value = "  +0.49%  "

$ws.Range("D40").Value = "'0.2326"
$ws.Range("B40").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4122) | Out-Null
$ws.Range("E40").Value = "  +3.26%  "

$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("D42").Value = "'0.6892"
$ws.Range("B42").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4122) | Out-Null
$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("D43").Value = "'1.259"
$ws.Range("B43").Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4122) | Out-Null
$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("D44").Value = "'14.87"
$ws.Range("B44").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4122) | Out-Null
$ws.Range("E44").Value = "  +6.22%  "

$ws.Range("D45").Value = "'0.6463"
$ws.Range("B45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4122) | Out-Null
$ws.Range("E45").Value = "  +1.93%  "

$ws.Range("D46").Value = "'2.324"
$ws.Range("B46").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4122) | Out-Null
$ws.Range("E46").Value = "  +5.79%  "

$ws.Range("D47").Value = "'0.00000000364"
$ws.Range("B47").Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4122) | Out-Null
$ws.Range("E47").Value = "  +21.91%  "

$ws.Range("D48").Value = "'3.690"
$ws.Range("B48").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4122) | Out-Null
$ws.Range("E48").Value = "  +1.76%  "

$ws.Range("D49").Value = "'1.257"
$ws.Range("B49").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4122) | Out-Null
$ws.Range("E49").Value = "  +1.07%  "

$ws.Range("D50").Value = "'83.57"
$ws.Range("B50").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4122) | Out-Null
$ws.Range("E50").Value = "  +2.18%  "

$ws.Range("D51").Value = "'0.3343"
$ws.Range("B51").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4122) | Out-Null
$ws.Range("E51").Value = "  +11.39%  "
